# Apply "testing demand increase through excel input" changes.
$wb = $excel.ActiveWorkbook

# --- Sheet "times": push StartTime/StopTime forward ---
$wsTimes = $wb.Worksheets.Item("times")
$wsTimes.Range("B2").Value = 47118.99861111111   # StartTime
$wsTimes.Range("B3").Value = 47482.99861111111   # StopTime

# --- Sheet "scenario_data_emlab": update scenario year and prices ---
$wsScenario = $wb.Worksheets.Item("scenario_data_emlab")
$wsScenario.Range("B1").Value = 2029             # scenario year
$wsScenario.Range("B2").Value = 50.28            # Co2Prices
$wsScenario.Range("B5").Value = 10.971           # FuelPrice_HARD_COAL
$wsScenario.Range("B6").Value = 26.934           # FuelPrice_NATURAL_GAS
$wsScenario.Range("B7").Value = 63.996           # FuelPrice_OIL
